$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.522.34'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.856.11'
$ws.Range('E3').Value = '  +1.27%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.71'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4747'
$ws.Range('E7').Value = '  +3.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2744'
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06318'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.74'
$ws.Range('E10').Value = '  +12.27%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.854.30'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07455'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.980'
$ws.Range('E13').Value = '  +2.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '84.66'
$ws.Range('E14').Value = '  +2.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6271'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').Value = '30.483.52'
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '246.52'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9999'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.72'
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007333'
$ws.Range('E20').Value = '  +2.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.945'
$ws.Range('E22').Value = '  +2.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.929'
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.138'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.45'
$ws.Range('E25').Value = '  -1.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.98'
$ws.Range('E26').Value = '  +3.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.876'
$ws.Range('E27').Value = '  +2.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1025'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  -1.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.025'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.838'
$ws.Range('E31').Value = '  +3.59%  '
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.135'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7012'
$ws.Range('E34').Value = '  +2.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.696'
$ws.Range('E35').Value = '  -0.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01899'
$ws.Range('E36').Value = '  +5.88%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.684'
$ws.Range('E37').Value = '  +3.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.001'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '106.60'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9999'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.539'
$ws.Range('E42').Value = '  +2.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4060'
$ws.Range('E43').Value = '  +3.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.228'
$ws.Range('E44').Value = '  +6.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.94'
$ws.Range('E45').Value = '  +7.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1203'
$ws.Range('E46').Value = '  +3.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '33.66'
$ws.Range('E47').Value = '  +4.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.569'
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3692'
$ws.Range('E51').Value = '  +3.35%  '
